# Weekly fruit/vegetable price update:
# Insert a new pair of rows (Primera/Segunda quality) for the new reporting
# date (serial 45021) right before the existing block of data that starts at
# row 985, pushing all subsequent rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 985:986 - this shifts the former rows 985-1056
# down to 987-1058 and keeps their values/formats intact.
$ws.Rows("985:986").Insert()

# --- New row 985 (Calidad = Primera) -------------------------------------
$ws.Range("A985").Value = 8
$ws.Range("B985").Value = "Terminal La Palmera de La Serena"
$ws.Range("C985").Value = "Coquimbo"
$ws.Range("D985").Value = 45021
$ws.Range("E985").Value = 4
$ws.Range("F985").Value = 100112023
$ws.Range("G985").Value = "Brócoli"
$ws.Range("H985").Value = "Sin especificar"
$ws.Range("I985").Value = "Primera"
$ws.Range("J985").Value = 2600
$ws.Range("K985").Value = 800
$ws.Range("L985").Value = 900
$ws.Range("M985").Value = 850
$ws.Range("N985").Value = "`$/unidad"
$ws.Range("O985").Value = "Provincia del Elquí"
$ws.Range("P985").Value = 850
$ws.Range("Q985").Value = 1
$ws.Range("R985").Value = "Hortaliza"

# --- New row 986 (Calidad = Segunda) --------------------------------------
$ws.Range("A986").Value = 8
$ws.Range("B986").Value = "Terminal La Palmera de La Serena"
$ws.Range("C986").Value = "Coquimbo"
$ws.Range("D986").Value = 45021
$ws.Range("E986").Value = 4
$ws.Range("F986").Value = 100112023
$ws.Range("G986").Value = "Brócoli"
$ws.Range("H986").Value = "Sin especificar"
$ws.Range("I986").Value = "Segunda"
$ws.Range("J986").Value = 1600
$ws.Range("K986").Value = 600
$ws.Range("L986").Value = 700
$ws.Range("M986").Value = 650
$ws.Range("N986").Value = "`$/unidad"
$ws.Range("O986").Value = "Provincia del Elquí"
$ws.Range("P986").Value = 650
$ws.Range("Q986").Value = 1
$ws.Range("R986").Value = "Hortaliza"

# Make sure the date cells keep the same date/time number format used by the
# rest of column D.
$ws.Range("D985:D986").NumberFormat = $ws.Range("D984").NumberFormat()
